{"js": "// Remove the TODO list item:\n// \"H. Add projectiles that rotate to face their target\"\n// while leaving the surrounding paragraphs untouched.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"Add projectiles that rotate to face their target\";\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.indexOf(targetText) !== -1) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nforeach ($p in @($d.Paragraphs)) {\n    if ($p.Range.Text -like \"*Add projectiles that rotate to face their target*\") {\n        $p.Range.Delete()\n    }\n}\n"}
